$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

$ws.Range('D2').Value = '28.222.89'
$ws.Range('E2').Value = '  -2.65%  '
$ws.Range('D3').Value = '1.867.89'
$ws.Range('E3').Value = '  -2.13%  '
$ws.Range('E4').Value = '  +0.18%  '
Set-TextValue $ws 'D5' '319.62'
$ws.Range('E5').Value = '  -1.40%  '
$ws.Range('E6').Value = '  +0.17%  '
Set-TextValue $ws 'D7' '0.4394'
$ws.Range('E7').Value = '  -4.28%  '
Set-TextValue $ws 'D8' '0.3687'
$ws.Range('E8').Value = '  -3.54%  '
Set-TextValue $ws 'D9' '0.07483'
$ws.Range('E9').Value = '  -2.96%  '
Set-TextValue $ws 'D10' '0.9347'
Set-TextValue $ws 'D11' '21.30'
$ws.Range('E11').Value = '  -3.41%  '
$ws.Range('D12').Value = '1.881.28'
$ws.Range('E12').Value = '  -0.44%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws 'D13' '5.468'
$ws.Range('E13').Value = '  -3.60%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws 'D14' '6.687'
$ws.Range('E14').Value = '  -3.53%  '
Set-TextValue $ws 'D15' '0.06919'
$ws.Range('E15').Value = '  -1.43%  '
Set-TextValue $ws 'D16' '1.004'
$ws.Range('E16').Value = '  +0.02%  '
Set-TextValue $ws 'D17' '82.07'
$ws.Range('E17').Value = '  -2.13%  '
Set-TextValue $ws 'D18' '0.000009005'
$ws.Range('E18').Value = '  -4.74%  '
$ws.Range('E19').Value = '  +0.22%  '
Set-TextValue $ws 'D20' '15.91'
$ws.Range('E20').Value = '  -4.71%  '
$ws.Range('D21').Value = '28.216.07'
$ws.Range('E21').Value = '  -2.55%  '
Set-TextValue $ws 'D22' '5.119'
$ws.Range('E22').Value = '  -3.74%  '
Set-TextValue $ws 'D23' '10.79'
$ws.Range('E23').Value = '  -0.82%  '
$ws.Range('D24').Value = '2.118.48'
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('E25').Value = '  -3.14%  '
$ws.Range('E26').Value = '  -2.12%  '
Set-TextValue $ws 'D27' '18.40'
Set-TextValue $ws 'D28' '5.321'
$ws.Range('E28').Value = '  -6.06%  '
Set-TextValue $ws 'D29' '113.42'
$ws.Range('E29').Value = '  -3.44%  '
Set-TextValue $ws 'D30' '1.725'
$ws.Range('E30').Value = '  -6.75%  '
Set-TextValue $ws 'D31' '0.09009'
$ws.Range('E31').Value = '  -2.79%  '
Set-TextValue $ws 'D32' '0.7935'
$ws.Range('E32').Value = '  -8.17%  '
Set-TextValue $ws 'D33' '4.844'
$ws.Range('E33').Value = '  -4.38%  '
Set-TextValue $ws 'D34' '1.172'
$ws.Range('E34').Value = '  -6.01%  '
Set-TextValue $ws 'D35' '2.923'
$ws.Range('E35').Value = '  -3.41%  '
$ws.Range('E36').Value = '  +0.20%  '
Set-TextValue $ws 'D37' '1.124'
$ws.Range('E37').Value = '  -2.71%  '
Set-TextValue $ws 'D38' '0.05444'
$ws.Range('E38').Value = '  -5.18%  '
Set-TextValue $ws 'D39' '0.01965'
$ws.Range('E39').Value = '  -3.76%  '
Set-TextValue $ws 'D40' '2.949'
$ws.Range('E40').Value = '  +3.46%  '
Set-TextValue $ws 'D41' '0.5253'
$ws.Range('E41').Value = '  -4.59%  '
Set-TextValue $ws 'D42' '7.045'
$ws.Range('E42').Value = '  -4.90%  '
Set-TextValue $ws 'D43' '0.1681'
$ws.Range('E43').Value = '  -4.33%  '
Set-TextValue $ws 'D44' '8.710'
$ws.Range('E44').Value = '  -6.50%  '
Set-TextValue $ws 'D45' '0.06743'
$ws.Range('E45').Value = '  -1.37%  '
Set-TextValue $ws 'D46' '0.4866'
$ws.Range('E46').Value = '  -6.03%  '
Set-TextValue $ws 'D47' '10.57'
$ws.Range('E47').Value = '  -6.58%  '
Set-TextValue $ws 'D48' '106.84'
$ws.Range('E48').Value = '  -3.69%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 'D49' '1.922'
$ws.Range('E49').Value = '  -6.04%  '
$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws 'D50' '1.002'
$ws.Range('E50').Value = '  +0.11%  '
Set-TextValue $ws 'D51' '1.671'
$ws.Range('E51').Value = '  -6.02%  '
